# Generate Report for Archive
#
# - Status text "Ready for handoff" -> "In Translation" everywhere it is used
#   (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share the same string).
# - The "Status" columns narrow to match the new (shorter) text:
#     Overview columns E & F, zh-cn column C, de-de column C.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update every cell carrying the old status text so the shared string is
# fully replaced (Overview has two status columns - zh-cn/de-de - each with
# two data rows; the other two sheets have one status column each).
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# Narrow the Status columns to fit the shorter text.
$newWidth = 12.576851254417766

$ws1.Columns.Item(5).ColumnWidth = $newWidth
$ws1.Columns.Item(6).ColumnWidth = $newWidth
$ws2.Columns.Item(3).ColumnWidth = $newWidth
$ws3.Columns.Item(3).ColumnWidth = $newWidth
